# add phoneNumber format control
# New rows of user-log data were appended to the "2022_2" sheet (rows 21-23).
# One of the new rows has a phoneNumber ("0022892942601") with leading zeros
# and no "+" prefix, which Excel would otherwise silently coerce to a plain
# number (stripping the leading zeros). To keep it intact we force column C
# (phoneNumber) to a Text number format before writing the value.
# The same text-format trick is used for column J (modelId), which this
# sheet always stores as text even when the value looks numeric (see the
# existing rows 9-20), while column I (userId) is kept as a real number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2022_2")

# Row 21
$ws.Cells.Item(21, 1).Value = "Tue Feb 22 2022"
$ws.Cells.Item(21, 2).Value = "17:14:42 GMT+0000 (Greenwich Mean Time)"
# phoneNumber (column C): this new entry has leading zeros and no "+"
# prefix, so Excel would otherwise auto-convert it to a number and lose
# the leading zeros. Force the cell to Text format before writing it.
$ws.Cells.Item(21, 3).NumberFormat = "@"
$ws.Cells.Item(21, 3).Value = "0022892942601"
$ws.Cells.Item(21, 4).Value = "User"
$ws.Cells.Item(21, 5).Value = "/api/auth/send-otp"
$ws.Cells.Item(21, 6).Value = "request"
$ws.Cells.Item(21, 7).Value = "failed"
$ws.Cells.Item(21, 8).Value = "0022892942601 request to receive otp"
$ws.Cells.Item(21, 11).Value = "error.invalid"

# Row 22
$ws.Cells.Item(22, 1).Value = "Tue Feb 22 2022"
$ws.Cells.Item(22, 2).Value = "17:18:25 GMT+0000 (Greenwich Mean Time)"
$ws.Cells.Item(22, 4).Value = "User"
$ws.Cells.Item(22, 5).Value = "/api/auth/login"
$ws.Cells.Item(22, 6).Value = "login"
$ws.Cells.Item(22, 7).Value = "failed"
$ws.Cells.Item(22, 8).Value = "invalid email  login"
$ws.Cells.Item(22, 11).Value = "error.invalid"
$ws.Cells.Item(22, 14).Value = "invalid email"

# Row 23
$ws.Cells.Item(23, 1).Value = "Tue Feb 22 2022"
$ws.Cells.Item(23, 2).Value = "17:20:47 GMT+0000 (Greenwich Mean Time)"
$ws.Cells.Item(23, 4).Value = "User"
$ws.Cells.Item(23, 5).Value = "/api/user/:id"
$ws.Cells.Item(23, 6).Value = "read"
$ws.Cells.Item(23, 7).Value = "succeeded"
$ws.Cells.Item(23, 8).Value = "Doe  John  read user 1"
$ws.Cells.Item(23, 9).Value = 1
# modelId (column J) is stored as text in this sheet even though it looks
# numeric (see existing rows 9-20), so keep it text here too.
$ws.Cells.Item(23, 10).NumberFormat = "@"
$ws.Cells.Item(23, 10).Value = "1"
$ws.Cells.Item(23, 12).Value = "Doe"
$ws.Cells.Item(23, 13).Value = "John"
